# Adds a new "2022-Q1" sheet (holding the fund-level breakdown, in the same
# layout as the existing quarterly sheets) right after "2021-Q4", and
# refreshes the "总计" (totals) rollup sheet with a new leading row for
# 2022-Q1, shifting the existing history down.

$wb = $excel.ActiveWorkbook

function Set-TextValue($rng, $val) {
    # Force the cell to literal text (matches fund codes like "002367"
    # keeping their leading zero, and numeric-looking ratios staying text)
    # while not leaving a lingering custom number-format style behind.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$q4Sheet = $wb.Worksheets.Item(5)

# --- 1. Insert the new "2022-Q1" sheet right after "2021-Q4" -------------
$newSheet = $wb.Worksheets.Add($null, $q4Sheet)
$newSheet.Name = "2022-Q1"

# Re-resolve the "总计" sheet AFTER the insert above: the sheet collection
# shifted, so any reference grabbed beforehand would now point at the
# wrong tab.
$totalSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Clone the "2021-Q4" layout/formatting (header row + column A index style)
# into the new sheet, then overwrite with the 2022-Q1 fund data.
$newSheet.Cells.Clear()
$q4Sheet.Range("A1:H5").Copy($newSheet.Range("A1"))
$newSheet.Range("A1").ClearContents()

Set-TextValue $newSheet.Range("B2") "257010"
$newSheet.Range("C2").Value = "国联安小盘精选混合"
Set-TextValue $newSheet.Range("D2") "9.15"
Set-TextValue $newSheet.Range("E2") "74.36"
Set-TextValue $newSheet.Range("F2") "4.55"
Set-TextValue $newSheet.Range("G2") "0.4163"
$newSheet.Range("H2").Value = 5

Set-TextValue $newSheet.Range("B3") "002367"
$newSheet.Range("C3").Value = "国联安安稳灵活配置混合"
Set-TextValue $newSheet.Range("D3") "2.32"
Set-TextValue $newSheet.Range("E3") "33.99"
Set-TextValue $newSheet.Range("F3") "3.24"
Set-TextValue $newSheet.Range("G3") "0.0752"
$newSheet.Range("H3").Value = 3

Set-TextValue $newSheet.Range("B4") "010746"
$newSheet.Range("C4").Value = "富安达长三角区域主题混合"
Set-TextValue $newSheet.Range("D4") "1.16"
Set-TextValue $newSheet.Range("E4") "92.45"
Set-TextValue $newSheet.Range("F4") "4.56"
Set-TextValue $newSheet.Range("G4") "0.0529"
$newSheet.Range("H4").Value = 9

Set-TextValue $newSheet.Range("B5") "006138"
$newSheet.Range("C5").Value = "国联安价值优选股票"
Set-TextValue $newSheet.Range("D5") "0.60"
Set-TextValue $newSheet.Range("E5") "93.30"
Set-TextValue $newSheet.Range("F5") "3.99"
Set-TextValue $newSheet.Range("G5") "0.0239"
$newSheet.Range("H5").Value = 9

# --- 2. Prepend a 2022-Q1 row to the "总计" rollup sheet ------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").Style = "Normal"
$totalSheet.Range("A3").Copy($totalSheet.Range("A2"))
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.57

# Column A is just a running 0-based row counter - renumber rows 3..7
# (originally 0..4) to 1..5 now that row 2 pushed them down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5

Write-Host "done"
